$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A; existing columns A:W shift to B:X
$ws.Columns("A").Insert()

# New "Match ID" header in row 2 (the visible header row)
$ws.Range("A2").Value = "Match ID"

# Apply bold font (no border) to the new column's header/data cells A2:A18
$ws.Range("A2:A18").Font.Bold = $true

# Fill in the Match ID value (12) for every data row
for ($r = 4; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value = 12
}

# Hidden totals row also gets the Match ID value, but keeps default (unstyled) formatting
$ws.Range("A19").Value = 12
$ws.Rows(19).AutoFit()

# Restore the selection to reflect the newly entered column of data
[void]$ws.Range("A2:A18").Select()
